$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the summary sheet ("总计"): insert a new first data row for
#    2022-Q4 and shift the existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

for ($r = 6; $r -ge 2; $r--) {
    $dstRow = $r + 1
    $summary.Range("B$dstRow").Value2 = $summary.Range("B$r").Value2
    $summary.Range("C$dstRow").Value2 = $summary.Range("C$r").Value2
    $summary.Range("D$dstRow").Value2 = $summary.Range("D$r").Value2
}

# New row 7 needs the same formatting as the rest of the index column;
# copy it (format + value) from A2, which already carries the right style.
$summary.Range("A2").Copy()
$summary.Range("A7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Renumber the 0-based index column for all 6 rows.
for ($r = 2; $r -le 7; $r++) {
    $summary.Range("A$r").Value2 = $r - 2
}

# Write the new 2022-Q4 summary row.
$summary.Range("B2").Value2 = "2022-Q4"
$summary.Range("C2").Value2 = 5
$summary.Range("D2").Value2 = 0.32

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet. Duplicating "2022-Q3" keeps every
#    formatting detail (borders, bold header, column widths, ...)
#    identical, and it is inserted right before "2022-Q3" — i.e. right
#    after "总计" — which matches the target tab order.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# The new sheet needs a 6th data row (Q3 only had 5); copy the format of
# row 5 down into row 6 first.
$newSheet.Range("A5").Copy($newSheet.Range("A6"))

# Fund-holding rows for 2022-Q4.
$data = @(
    @(0, "005585", "银河文体娱乐主题灵活配置混合A", "3.15", "88.90", "4.94", "0.1556", 8),
    @(1, "001628", "招商体育文化休闲股票A",         "2.33", "93.03", "4.92", "0.1146", 7),
    @(2, "003397", "银华体育文化灵活配置混合",       "0.53", "87.74", "3.78", "0.0200", 6),
    @(3, "015667", "银河文体娱乐主题灵活配置混合C", "0.38", "88.90", "4.94", "0.0188", 8),
    @(4, "015395", "招商体育文化休闲股票C",         "0.29", "93.03", "4.92", "0.0143", 7)
)

# B, D, E, F, G hold numeric-looking text (fund code / percentages / nav)
# that must stay text (leading zeros, trailing zeros) — same as every
# other quarter sheet in this workbook. Mark them as Text before writing,
# then drop back to the default "Normal" style so no stray formatting is
# left behind (only the value's string-ness is kept).
$textCols = @("B", "D", "E", "F", "G")

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]

    $newSheet.Range("A$row").Value2 = $rowData[0]

    foreach ($col in $textCols) {
        $newSheet.Range("$col$row").NumberFormat = "@"
    }
    $newSheet.Range("B$row").Value2 = $rowData[1]
    $newSheet.Range("C$row").Value2 = $rowData[2]
    $newSheet.Range("D$row").Value2 = $rowData[3]
    $newSheet.Range("E$row").Value2 = $rowData[4]
    $newSheet.Range("F$row").Value2 = $rowData[5]
    $newSheet.Range("G$row").Value2 = $rowData[6]
    foreach ($col in $textCols) {
        $newSheet.Range("$col$row").Style = "Normal"
    }

    $newSheet.Range("H$row").Value2 = $rowData[7]
}
